$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

$title = $s1.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "."
$titleRange.Text = "A Table, with a caption"

$caption = $s1.Shapes.Item(3)
$captionRange = $caption.TextFrame.TextRange
$captionRange.Text = "."
$captionRange.Text = "Demonstration of simple table syntax, with alignment"
